$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overall" (sheet2) - row 2 values updated
# ---------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")

$wsOverall.Range("B2").Value = 36
$wsOverall.Range("C2").Value = 16
$wsOverall.Range("D2").Value = 0.6527948839311194
$wsOverall.Range("E2").Value = 0.51851851851851849
$wsOverall.Range("F2").Value = 0.8087287276360754
$wsOverall.Range("G2").Value = 25
$wsOverall.Range("H2").Value = 5
$wsOverall.Range("I2").Value = 30
$wsOverall.Range("J2").Value = 346
$wsOverall.Range("K2").Value = 11

# ---------------------------------------------------------------
# Sheet "Zones" (sheet3) - rows 2..14 values updated
# ---------------------------------------------------------------
$wsZones = $wb.Worksheets.Item("Zones")

# row 2 - Zone 1
$wsZones.Range("B2").Value = 3
$wsZones.Range("C2").Value = 1
$wsZones.Range("D2").Value = 0.65999999999999992
$wsZones.Range("E2").Value = 0.4444444444444442
$wsZones.Range("F2").Value = 0.98333333333333361

# row 3 - Zone 2
$wsZones.Range("B3").Value = 3
$wsZones.Range("C3").Value = 1
$wsZones.Range("D3").Value = 0.71875
$wsZones.Range("E3").Value = 0.82777777777777783
$wsZones.Range("F3").Value = 0.65333333333333321

# row 4 - Zone 3
$wsZones.Range("B4").Value = 2
$wsZones.Range("C4").Value = 2
$wsZones.Range("D4").Value = 0.52619047619047621
$wsZones.Range("E4").Value = 0.68333333333333313
$wsZones.Range("F4").Value = 0.46333333333333349

# row 5 - Zone 4
$wsZones.Range("B5").Value = 2
$wsZones.Range("C5").Value = 1
$wsZones.Range("D5").Value = 0.80833333333333324
$wsZones.Range("E5").Value = 1.0249999999999999
$wsZones.Range("F5").Value = 0.59166666666666656

# row 6 - Zone 5 (F6 removed, E6 now equals D6)
$wsZones.Range("B6").Value = 11
$wsZones.Range("C6").Value = 0
$wsZones.Range("D6").Value = 0.41212121212121217
$wsZones.Range("E6").Value = 0.41212121212121217
$wsZones.Range("F6").ClearContents()

# row 7 - Zone 6
$wsZones.Range("B7").Value = 2
$wsZones.Range("C7").Value = 0
$wsZones.Range("D7").Value = 0.54999999999999993
$wsZones.Range("E7").Value = 0.24999999999999994
$wsZones.Range("F7").Value = 1.1499999999999999

# row 8 - Zone 7
$wsZones.Range("B8").Value = 2
$wsZones.Range("C8").Value = 1
$wsZones.Range("D8").Value = 0.48333333333333323
$wsZones.Range("E8").Value = 0.44166666666666649
$wsZones.Range("F8").Value = 0.56666666666666687

# row 9 - Zone 8
$wsZones.Range("B9").Value = 5
$wsZones.Range("C9").Value = 1
$wsZones.Range("D9").Value = 0.71111111111111103
$wsZones.Range("E9").Value = 0.65666666666666662
$wsZones.Range("F9").Value = 0.77916666666666656

# row 10 - Zone 9
$wsZones.Range("B10").Value = 0
$wsZones.Range("C10").Value = 2
$wsZones.Range("D10").Value = 0.73333333333333328
$wsZones.Range("F10").Value = 0.73333333333333328

# row 11 - Zone 10
$wsZones.Range("B11").Value = 2
$wsZones.Range("C11").Value = 2
$wsZones.Range("D11").Value = 0.8
$wsZones.Range("E11").Value = 0.25000000000000011
$wsZones.Range("F11").Value = 1.1666666666666667

# row 12 - Zone 11 (F12 removed, E12 now equals D12)
$wsZones.Range("B12").Value = 3
$wsZones.Range("C12").Value = 0
$wsZones.Range("D12").Value = 0.37222222222222218
$wsZones.Range("E12").Value = 0.37222222222222218
$wsZones.Range("F12").ClearContents()

# row 13 - Zone 12 (E13 removed, D13 and F13 now equal)
$wsZones.Range("B13").Value = 0
$wsZones.Range("C13").Value = 5
$wsZones.Range("D13").Value = 1.1138888888888892
$wsZones.Range("E13").ClearContents()
$wsZones.Range("F13").Value = 1.1138888888888892

# row 14 - Zone 13 (F14 removed, E14 added)
$wsZones.Range("B14").Value = 1
$wsZones.Range("C14").Value = 0
$wsZones.Range("D14").Value = 0.6166666666666667
$wsZones.Range("E14").Value = 0.6166666666666667
$wsZones.Range("F14").ClearContents()
